$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 was a blank template row; turn it into a real client row ("NILTON
# BISPO") by first cloning the formatting of the row above it (row 12,
# which already has the "filled" look: borders/fill/font + the date and
# integer number formats used by the rest of the table) and then filling
# in the values.
$ws.Range("B12:E12").Copy()
$ws.Range("B13:E13").PasteSpecial(-4122)

$ws.Range("B13").Value = "NILTON BISPO"
$ws.Range("C13").Value = "7fb8603912d23d36b1c4c22f51f81385"
$ws.Range("D13").Value = 44830
$ws.Range("E13").Value = 9

# Row 16 ("TESTE EDSON") has its activation code (column C) cleared out.
$ws.Range("C16").Value = ""
